# Fix for new web page
# Adds three new header columns to the company list sheet:
#   - "聯絡手機" (Contact Mobile) before "聯絡電話"
#   - "郵遞區號" (Postal Code) before "通訊地址"
#   - "相關資料" (Related Info) appended after the last column ("公司網址")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column before "通訊地址" (currently column I) and label it.
$ws.Columns("I:I").Insert()
$ws.Range("I1").Value = "郵遞區號"

# 2) Append a new column after the current last column ("公司網址", now column O)
#    and label it.
$ws.Range("P1").Value = "相關資料"

# 3) Insert a new column before "聯絡電話" (still column F) and label it.
$ws.Columns("F:F").Insert()
$ws.Range("F1").Value = "聯絡手機"

# Restore a sensible active selection like the edited workbook.
$ws.Range("H4").Select()
